# Apply weekly fruit/vegetable price update: rows 2-12 get their
# D (Fecha), K (Variedad), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), Q (Unidad de
# comercializacion), R (Origen), S (Precio $/Kg) and T (Kg / unidad)
# values reshuffled across the rows, per the new weekly data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values keyed by destination row number.
$data = @{
    2  = @{ D = 44210; K = "Rainier";     L = "Segunda"; M = 250; N = 21000; O = 22000; P = 21500; Q = "`$/caja 18 kilos";    R = "Región de O'Higgins"; S = 1194; T = 18 }
    3  = @{ D = 44568; K = "Santina";     L = "Segunda"; M = 200; N = 15000; O = 16000; P = 15500; Q = "`$/bandeja 12 kilos"; R = "Región de O'Higgins"; S = 1292; T = 12 }
    4  = @{ D = 44571; K = "Brooks";      L = "Segunda"; M = 400; N = 8500;  O = 9000;  P = 8750;  Q = "`$/bandeja 10 kilos"; R = "Región de O'Higgins"; S = 875;  T = 10 }
    5  = @{ D = 44229; K = "Santina";     L = "Primera"; M = 250; N = 6500;  O = 7000;  P = 6750;  Q = "`$/bandeja 5 kilos";  R = "Provincia de Curicó"; S = 1350; T = 5  }
    6  = @{ D = 44208; K = "Lapins";      L = "Segunda"; M = 200; N = 10500; O = 11000; P = 10750; Q = "`$/bandeja 12 kilos"; R = "Provincia de Curicó"; S = 896;  T = 12 }
    7  = @{ D = 44537; K = "Brooks";      L = "Primera"; M = 200; N = 29000; O = 30000; P = 29500; Q = "`$/caja 20 kilos";    R = "Región de O'Higgins"; S = 1475; T = 20 }
    8  = @{ D = 44557; K = "Lapins";      L = "Primera"; M = 250; N = 9000;  O = 10000; P = 9500;  Q = "`$/bandeja 10 kilos"; R = "Provincia de Curicó"; S = 950;  T = 10 }
    9  = @{ D = 44580; K = "Sweet Heart"; L = "Segunda"; M = 300; N = 7000;  O = 8000;  P = 7500;  Q = "`$/bandeja 10 kilos"; R = "Región de O'Higgins"; S = 750;  T = 10 }
    10 = @{ D = 44175; K = "Rainier";     L = "Segunda"; M = 270; N = 25000; O = 26000; P = 25500; Q = "`$/caja 18 kilos";    R = "Región de O'Higgins"; S = 1417; T = 18 }
    11 = @{ D = 44161; K = "Bing";        L = "Primera"; M = 160; N = 39000; O = 40000; P = 39500; Q = "`$/caja 20 kilos";    R = "Provincia de Curicó"; S = 1975; T = 20 }
    12 = @{ D = 44532; K = "Brooks";      L = "Primera"; M = 400; N = 27000; O = 28000; P = 27500; Q = "`$/bandeja 12 kilos"; R = "Región de O'Higgins"; S = 2292; T = 12 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $vals.D   # D: Fecha
    $ws.Cells.Item($row, 11).Value = $vals.K   # K: Variedad
    $ws.Cells.Item($row, 12).Value = $vals.L   # L: Calidad
    $ws.Cells.Item($row, 13).Value = $vals.M   # M: Volumen
    $ws.Cells.Item($row, 14).Value = $vals.N   # N: Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals.O   # O: Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals.P   # P: Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = $vals.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($row, 18).Value = $vals.R   # R: Origen
    $ws.Cells.Item($row, 19).Value = $vals.S   # S: Precio $/Kg
    $ws.Cells.Item($row, 20).Value = $vals.T   # T: Kg / unidad
}
